$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookups")

# The "Gender Not Used" (missing) entry moves from the Gender lookup
# column (F/G) to the GenderB lookup column (H/I), shifting the
# remaining GenderB entries down by one row.

$ws.Range("F2").Value = "male"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = "Gender Not Used"
$ws.Range("I2").Value = 0

$ws.Range("F3").Value = "female"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = "male"
$ws.Range("I3").Value = 1

$ws.Range("F4").Value = $null
$ws.Range("G4").Value = $null
$ws.Range("H4").Value = "female"
$ws.Range("I4").Value = 2
